$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.665.45'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '2.044.24'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''227.41'
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").Value = '''59.28'
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -2.06%  '
$ws.Range("E10").Value = '  +3.06%  '
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").Value = '2.346.35'
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("D13").Value = '''14.38'
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").Value = '''21.33'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = '''5.46'
$ws.Range("E15").Value = '  +5.38%  '
$ws.Range("D16").Value = '''0.761'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").Value = '2.043.02'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '37.678.99'
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").Value = '''69.37'
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("D20").Value = '''5.91'
$ws.Range("E20").Value = '  -1.82%  '
$ws.Range("D21").Value = '0.0₃0830'
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("D22").Value = '''222.62'
$ws.Range("E22").Value = '  -1.08%  '
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").Value = '''2.28'
$ws.Range("E25").Value = '  +2.81%  '
$ws.Range("D26").Value = '''168.99'
$ws.Range("E26").Value = '  +2.34%  '
$ws.Range("D27").Value = '''9.29'
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("D29").Value = '''18.75'
$ws.Range("E29").Value = '  -1.15%  '
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("E31").Value = '  -0.65%  '
$ws.Range("E32").Value = '  +8.41%  '
$ws.Range("D33").Value = '''4.37'
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("D34").Value = '''4.54'
$ws.Range("E34").Value = '  +0.61%  '
$ws.Range("D35").Value = '''0.0602'
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = '''6.47'
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("E37").Value = '  +3.59%  '
$ws.Range("D38").Value = '''3.45'
$ws.Range("E38").Value = '  +6.43%  '
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("D40").Value = '''18.37'
$ws.Range("E40").Value = '  +9.70%  '
$ws.Range("D41").Value = '1.527.24'
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").Value = '''97.92'
$ws.Range("E42").Value = '  +1.44%  '
$ws.Range("E43").Value = '  -1.98%  '
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").Value = '''4.17'
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("E46").Value = '  -2.89%  '
$ws.Range("D47").Value = '''1.11'
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("E48").Value = '  -0.90%  '
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").Value = '''7.06'
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").Value = '2.235.54'
$ws.Range("E51").Value = '  +0.78%  '
